# Team_Everyday_Attendence.xlsx - "Team attendance [B8-G1] 04-08"
#
# Adds a new team member column (Priya Gawhane, col I) and a new day's
# attendance row (04-08, row 4), plus the "No Response" follow-up
# comments RENUKA left on the still-unanswered cells, and moves the
# active selection to the last-edited cell (I4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header: I1 = "Priya Gawhane" ---------------------------
$ws.Range("I1").Value = "Priya Gawhane"

# --- Existing rows 2 & 3 get an ABSENT mark for the new column ---------
$ws.Range("I2").Value = "ABSENT"
$ws.Range("I3").Value = "ABSENT"

# --- New attendance row for 04-08 (row 4) -------------------------------
# A4 (the date, 45142) already exists in the sheet.
$ws.Range("B4").Value = "PRESENT"
$ws.Range("C4").Value = "PRESENT"
$ws.Range("D4").Value = "PRESENT"
$ws.Range("E4").Value = "PRESENT"
$ws.Range("F4").Value = "PRESENT"
$ws.Range("G4").Value = "ABSENT"
$ws.Range("H4").Value = "ABSENT"
$ws.Range("I4").Value = "ABSENT"

# --- RENUKA's "No Response" comments on the newly-absent cells ---------
$ws.Range("I2").AddComment("RENUKA:`nNo Response") | Out-Null
$ws.Range("I3").AddComment("RENUKA:`nNo Response") | Out-Null
$ws.Range("G4").AddComment("RENUKA:`nNo Response`n") | Out-Null
$ws.Range("H4").AddComment("RENUKA:`nNo Response") | Out-Null
$ws.Range("I4").AddComment("RENUKA:`nNo Response") | Out-Null

# --- Move the selection to the last cell that was filled in ------------
$ws.Range("I4").Select() | Out-Null
